$d = $word.ActiveDocument

# Move to the very end of the document content (after "Deliverable 2")
$rng = $d.Content
$rng.Collapse(0)  # wdCollapseEnd

$rng.InsertAfter([char]13 + [char]13 + "This text is modified by NMS")
